$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $val) {
    $r = $ws.Range($rangeAddr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "304.99"
Set-TextValue "E2" "-0.07%"

Set-TextValue "D3" "38.32"
Set-TextValue "E3" "7.60%"

Set-TextValue "D4" "5.098"
Set-TextValue "E4" "0.90%"

Set-TextValue "D5" "0.08047"
Set-TextValue "E5" "0.61%"

Set-TextValue "D6" "1.920"
Set-TextValue "E6" "3.06%"

$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue "D7" "4.199"
Set-TextValue "E7" "1.79%"

$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
Set-TextValue "D8" "7.946"
Set-TextValue "E8" "2.17%"

$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D9" "0.9278"
Set-TextValue "E9" "0.71%"

$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D10" "0.1440"
Set-TextValue "E10" "11.61%"

$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D11" "0.1921"
Set-TextValue "E11" "2.24%"

$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D12" "0.08935"
Set-TextValue "E12" "-0.82%"

$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D13" "0.03531"
Set-TextValue "E13" "3.00%"

$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D14" "0.09772"
Set-TextValue "E14" "-0.94%"

$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D15" "0.001398"
Set-TextValue "E15" "-0.38%"

$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D16" "0.005975"
Set-TextValue "E16" "-4.07%"

$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D17" "3.724"
Set-TextValue "E17" "-3.09%"

Set-TextValue "D18" "3.426"
Set-TextValue "E18" "1.45%"

Set-TextValue "E19" "1.10%"

Set-TextValue "D20" "0.1334"
Set-TextValue "E20" "-0.52%"

Set-TextValue "D21" "4.840"
Set-TextValue "E21" "1.06%"

Set-TextValue "D22" "0.2409"
Set-TextValue "E22" "4.31%"

Set-TextValue "D23" "0.04363"
Set-TextValue "E23" "-1.40%"

Set-TextValue "D24" "0.001228"
Set-TextValue "E24" "-0.49%"

Set-TextValue "D25" "0.004110"
Set-TextValue "E25" "-15.51%"

Set-TextValue "E27" "0.42%"

Set-TextValue "D39" "0.02066"
Set-TextValue "E39" "7.63%"

Set-TextValue "D40" "0.05018"
Set-TextValue "E40" "-2.80%"

Set-TextValue "D41" "0.007490"
Set-TextValue "E41" "-0.84%"

Set-TextValue "D42" "0.01013"
Set-TextValue "E42" "-0.24%"

Set-TextValue "E43" "0.13%"

Set-TextValue "D44" "0.002144"
Set-TextValue "E44" "0.90%"

Set-TextValue "D45" "0.008910"
Set-TextValue "E45" "-9.79%"

Set-TextValue "E46" "-0.12%"

Set-TextValue "E47" "0.12%"

Set-TextValue "D48" "0.002991"

Set-TextValue "D49" "0.001604"
Set-TextValue "E49" "28.09%"

Set-TextValue "D50" "0.00002106"
Set-TextValue "E50" "0.12%"

Set-TextValue "D51" "0.0002006"
Set-TextValue "E51" "0.12%"
